# stats_country_ranks: recalc "completeness_rank" (column E) now that the
# missing/non-zero daily-case stats feeding it were corrected. The rank is
# a dense, tie-broken ordering of column D (nazero_daily_cases); fixing the
# tie-break logic reshuffles the rank among rows that previously shared the
# same D value. Column E is independent data (no formula on the sheet), so
# the corrected ranks are written straight to each affected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$completenessRankFixes = @{
    31  = 128
    32  = 127
    44  = 114
    45  = 115
    48  = 110
    49  = 111
    51  = 107
    52  = 108
    53  = 105
    54  = 106
    69  = 89
    70  = 90
    76  = 83
    77  = 80
    78  = 82
    79  = 81
    80  = 78
    81  = 79
    82  = 76
    83  = 77
    86  = 74
    87  = 73
    88  = 71
    89  = 68
    92  = 70
    97  = 59
    98  = 60
    100 = 62
    101 = 58
    102 = 57
    104 = 55
    105 = 56
    110 = 48
    111 = 49
    112 = 47
    113 = 45
    114 = 46
    117 = 40
    118 = 39
    119 = 42
    120 = 38
    121 = 41
    122 = 35
    123 = 36
    124 = 37
    125 = 34
    127 = 31
    128 = 33
    129 = 27
    131 = 30
    133 = 25
    134 = 26
    135 = 23
    136 = 24
    138 = 21
    139 = 20
    140 = 18
    142 = 17
    143 = 16
    148 = 10
    149 = 11
    150 = 8
    151 = 9
    155 = 4
    156 = 3
}

foreach ($row in $completenessRankFixes.Keys) {
    $ws.Range("E$row").Value = $completenessRankFixes[$row]
}
